$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title row
$ws.Range("A1").Value = "Thống kê từ ngày 2025-05-23 đến 2026-06-01"

# Header row (unchanged, but re-assert for safety)
$ws.Range("A3").Value = "Tên phim"
$ws.Range("B3").Value = "Số suất chiếu"
$ws.Range("C3").Value = "Số vé bán"
$ws.Range("D3").Value = "Doanh thu (VNĐ)"

# Data rows 4-10
$ws.Range("A4").Value = "Mario Bros"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 13
$ws.Range("D4").Value = 750000

$ws.Range("A5").Value = "Rocky"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 280000

$ws.Range("A6").Value = "Ant-Man"
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

$ws.Range("A7").Value = "Avengers"
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0

$ws.Range("A8").Value = "Immaculate"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0

$ws.Range("A9").Value = "Ròm"
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0

# New row 10 for "The Lion King" (previously the totals row was row 11; now blank row moves to 11)
$ws.Range("A10").Value = "The Lion King"
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0

# Clear old totals row (11) leftover cells since totals now live on row 12
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""

# Totals row moved from 11 to 12
$ws.Range("C12").Value = "Tổng doanh thu:"
$ws.Range("D12").Value = 1030000
